$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$cell = $ws1.Range("E2")
$cell.Value = "domestic water usage, Vww (L Cap-1 Day-1)"
$cell.Characters(1, 23).Font.Bold = $true
$cell.Characters(24, 2).Font.Bold = $true
$cell.Characters(24, 2).Font.Subscript = $true
$cell.Characters(26, 17).Font.Bold = $true
